$d = $word.ActiveDocument

function Get-ParaIndexByText($pattern) {
    $count = $word.ActiveDocument.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $t = $word.ActiveDocument.Paragraphs.Item($i).Range.Text
        if ($t -match $pattern) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# Edit 1: "Empowering STEM members ... career opportunities." paragraph
#         -> replaced with "I serve on the Leadership Committee, ..." and a
#         new empty ListParagraph paragraph is inserted right after it.
# ---------------------------------------------------------------------------
$idx1 = Get-ParaIndexByText("Empowering STEM members")
$p1 = $word.ActiveDocument.Paragraphs.Item($idx1)

$xmlHeader = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">'
$xmlFooter = '</w:document></pkg:xmlData></pkg:part></pkg:package>'

$frag1 = $xmlHeader + '<w:body>' + `
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:ind w:right="0"/></w:pPr><w:r><w:t>I serve on the Leadership Committee, where I help shape the overall direction of KNSBE and support our Admin Zone.</w:t></w:r></w:p>' + `
    '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:ind w:right="0"/></w:pPr></w:p>' + `
    '</w:body>' + $xmlFooter

$p1.Range.InsertXML($frag1)

# The freshly-inserted empty paragraph needs an explicit firstLine="0" on its
# <w:ind>; assigning it through the paragraph-format object model (rather
# than the XML fragment) makes the zero value persist in the saved XML.
$emptyParaIdx = $idx1 + 1
$emptyPara = $word.ActiveDocument.Paragraphs.Item($emptyParaIdx)
$emptyPara.Range.ParagraphFormat.FirstLineIndent = 0

# ---------------------------------------------------------------------------
# Edit 2: "I am a part of the Leadership Committee: ..." paragraph
#         -> pStyle "ListParagraph" removed (keeps numPr/ind) and text
#         replaced with the new Service Desk description.
# ---------------------------------------------------------------------------
$idx2 = Get-ParaIndexByText("I am a part of the Leadership Committee")
$p2 = $word.ActiveDocument.Paragraphs.Item($idx2)

$frag2 = $xmlHeader + '<w:body>' + `
    '<w:p><w:pPr><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:ind w:right="0"/></w:pPr><w:r><w:t>Welcome guests, answer questions, manage memberships and equipment rentals, ensure safety. Gain customer service experience in campus recreation.</w:t></w:r></w:p>' + `
    '</w:body>' + $xmlFooter

$p2.Range.InsertXML($frag2)

# ---------------------------------------------------------------------------
# Edit 3: styles.xml - remove <w:semiHidden/> from the DefaultParagraphFont
#         character style.
# ---------------------------------------------------------------------------
$styles = $word.ActiveDocument.Styles
$dpf = $styles.Item("Default Paragraph Font")
$dpf.SemiHidden = $false

Write-Host "Edits applied successfully"
